$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = -3
$ws.Range("F5").Value = 1
$ws.Range("F8").Value = -3
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = -5
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = -3
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 0
$ws.Range("F22").Value = 3
$ws.Range("F25").Value = -2
$ws.Range("F27").Value = -2
$ws.Range("F31").Value = 3
$ws.Range("F37").Value = 0
$ws.Range("F39").Value = 2
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 3
$ws.Range("F46").Value = -1
$ws.Range("F48").Value = 1
$ws.Range("F50").Value = -1
$ws.Range("F55").Value = 2
$ws.Range("F57").Value = -3
$ws.Range("F58").Value = -1
$ws.Range("F60").Value = 1
$ws.Range("F62").Value = -2
$ws.Range("F67").Value = 6
$ws.Range("F68").Value = 0
